$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet tab from "1" to "გარდაბანი"
$ws.Name = "გარდაბანი"

# Clear the subtitle text in row 2 (A2), keep the (now-blank) row
$ws.Range("A2").ClearContents()

# Delete row 3 (was blank spacer row) entirely - rows shift up
$ws.Rows(3).Delete()

# Delete columns B:C (1989 / 2002 figures), leaving only the 2014 column,
# which shifts left into column B
$ws.Range("B:C").Delete()

# Update selection to A2, matching the saved view state
$ws.Range("A2").Select()

$wb.Save()
